$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number but must remain
# stored as literal text (matching the workbook's original inline-string
# representation, e.g. "256.80", "0.0360", "2.20"). Force text format first
# so Excel does not silently coerce these into numeric cells, then clear the
# number-format override afterwards so no stray cell style is left behind.
$textCoercedCells = @(
    "D5"
    "D6"
    "D7"
    "D10"
    "D11"
    "D12"
    "D13"
    "D17"
    "D19"
    "D22"
    "D23"
    "D24"
    "D26"
    "D27"
    "D31"
    "D32"
    "D33"
    "D36"
    "D37"
    "D38"
    "D39"
    "D40"
    "D44"
    "D45"
    "D47"
    "D48"
    "D49"
    "D50"
)
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price column cells that need text coercion ---
$ws.Range("D5").Value = "256.80"
$ws.Range("D6").Value = "0.615"
$ws.Range("D7").Value = "76.35"
$ws.Range("D10").Value = "41.91"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D12").Value = "6.93"
$ws.Range("D13").Value = "0.103"
$ws.Range("D17").Value = "0.785"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("D22").Value = "230.13"
$ws.Range("D23").Value = "2.20"
$ws.Range("D24").Value = "9.27"
$ws.Range("D26").Value = "42.23"
$ws.Range("D27").Value = "10.75"
$ws.Range("D31").Value = "174.14"
$ws.Range("D32").Value = "20.29"
$ws.Range("D33").Value = "0.0871"
$ws.Range("D36").Value = "0.0360"
$ws.Range("D37").Value = "0.106"
$ws.Range("D38").Value = "4.30"
$ws.Range("D39").Value = "12.75"
$ws.Range("D40").Value = "2.84"
$ws.Range("D44").Value = "60.03"
$ws.Range("D45").Value = "102.61"
$ws.Range("D47").Value = "8.30"
$ws.Range("D48").Value = "0.455"
$ws.Range("D49").Value = "1.11"
$ws.Range("D50").Value = "1.13"

foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).ClearFormats()
}

# --- Remaining cells (coin names, links, prices with multiple dots, percentages) ---
$ws.Range("D2").Value = "42.919.96"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.212.76"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +4.79%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "2.545.79"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "2.215.38"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "42.877.55"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -7.70%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +6.70%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -4.84%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E33").Value = "  +9.01%  "
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  +7.87%  "
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("E40").Value = "  +17.54%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("E48").Value = "  -6.39%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.436.47"
$ws.Range("E51").Value = "  -1.16%  "
